$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "userName String"
$ws.Range("B6").Value = "newUserName String"
$ws.Range("B7").Value = "newEmail"
$ws.Range("B8").Value = "newName"

$ws.Range("B23").Value = "userName = newUserName"
$ws.Range("C23").Value = "\\admin able to edit users information"
$ws.Range("B24").Value = "email = newEmail"
$ws.Range("B25").Value = "name - newName"

$ws.Range("B29").Value = "commit changes"

$ws.Range("B25").Select()
